$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column G width update (was 24.1328125 -> now ~36.3984375, bestFit/customWidth) ---
$ws.Columns.Item(7).ColumnWidth = 35.5

# --- Copy formatting (style) from column F into the new column G for the affected rows ---
$rows = @(14,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32)
foreach ($r in $rows) {
    $ws.Range("F$r").Copy($ws.Range("G$r"))
}

# --- Set the new/changed values in column G ---
$ws.Range("G19").Value = '$model.setDefaultReturnQueue($param);'
$ws.Range("G20").Value = "Default return queue"
$ws.Range("G21").Value = "null"
$ws.Range("G22").Value = "null"
$ws.Range("G23").Value = "null"
$ws.Range("G24").Value = "null"
$ws.Range("G25").Value = '"Fulfill"'
$ws.Range("G26").Value = '"Fulfill"'
$ws.Range("G27").Value = '"Fulfill"'
$ws.Range("G28").Value = '"Fulfill"'
$ws.Range("G29").Value = '"Fulfill"'
$ws.Range("G30").Value = '"Fulfill"'
$ws.Range("G31").Value = '"Fulfill"'
$ws.Range("G32").Value = "null"

# --- Update sheet view scroll position / selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 4
$ws.Range("G19").Select()
